# Updated cryptos list on Sun Apr 16 05:18:24 UTC 2023 with GitHub Actions
#
# Applies the latest scraped Price (column D) / Volume(1h) (column E) values
# to the cryptos table, plus the PancakeSwap <-> BabyDogeCoin row swap
# (rows 47/48) where name, link, price and volume all change together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    # Force the cell to be written as text even when the string looks like a
    # number (e.g. "1.009"), then drop back to the default "Normal" style so
    # no stray number-format/style index is left behind on the cell.
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "30.611.85"
Set-TextCell "E2" "  +0.60%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.117.22"
Set-TextCell "E3" "  +1.22%  "

# Row 4 - TetherUSD (price unchanged)
Set-TextCell "E4" "  +0.86%  "

# Row 5 - BNB
Set-TextCell "D5" "336.79"
Set-TextCell "E5" "  +2.01%  "

# Row 6 - USDC
Set-TextCell "D6" "1.009"
Set-TextCell "E6" "  +0.86%  "

# Row 7 - XRP
Set-TextCell "D7" "0.5242"
Set-TextCell "E7" "  +0.31%  "

# Row 8 - Cardano
Set-TextCell "D8" "0.4563"
Set-TextCell "E8" "  +4.33%  "

# Row 9 - OKB
Set-TextCell "D9" "54.63"
Set-TextCell "E9" "  +2.07%  "

# Row 10 - Dogecoin
Set-TextCell "D10" "0.09098"
Set-TextCell "E10" "  +2.40%  "

# Row 11 - Polygon
Set-TextCell "D11" "1.171"
Set-TextCell "E11" "  +1.49%  "

# Row 12 - Solana
Set-TextCell "D12" "24.62"
Set-TextCell "E12" "  +1.43%  "

# Row 13 - WrappedEther
Set-TextCell "D13" "2.111.47"
Set-TextCell "E13" "  +1.12%  "

# Row 14 - Polkadot
Set-TextCell "D14" "6.852"
Set-TextCell "E14" "  +2.23%  "

# Row 15 - Chainlink
Set-TextCell "D15" "8.098"
Set-TextCell "E15" "  +5.39%  "

# Row 16 - ShibaInu
Set-TextCell "D16" "0.00001176"
Set-TextCell "E16" "  +4.93%  "

# Row 17 - Litecoin
Set-TextCell "D17" "97.06"
Set-TextCell "E17" "  +1.29%  "

# Row 18 - BinanceUSD (price unchanged)
Set-TextCell "E18" "  +0.76%  "

# Row 19 - TRON
Set-TextCell "D19" "0.06687"
Set-TextCell "E19" "  +1.45%  "

# Row 20 - Avalanche
Set-TextCell "D20" "19.40"
Set-TextCell "E20" "  +0.97%  "

# Row 21 - Dai
Set-TextCell "D21" "1.009"
Set-TextCell "E21" "  +0.84%  "

# Row 22 - Uniswap
Set-TextCell "D22" "6.297"
Set-TextCell "E22" "  +0.47%  "

# Row 23 - WrappedBTC
Set-TextCell "D23" "30.683.16"
Set-TextCell "E23" "  +0.71%  "

# Row 24 - Cosmos
Set-TextCell "D24" "12.85"
Set-TextCell "E24" "  +4.98%  "

# Row 25 - Toncoin
Set-TextCell "D25" "2.352"
Set-TextCell "E25" "  +0.94%  "

# Row 26 - WrappedliquidstakedEther2.0
Set-TextCell "D26" "2.361.18"
Set-TextCell "E26" "  +1.18%  "

# Row 27 - EthereumClassic
Set-TextCell "D27" "22.39"
Set-TextCell "E27" "  +0.63%  "

# Row 28 - Monero
Set-TextCell "D28" "163.83"
Set-TextCell "E28" "  +0.59%  "

# Row 29 - LidoDAOToken
Set-TextCell "D29" "2.540"
Set-TextCell "E29" "  -1.01%  "

# Row 30 - BitcoinCash
Set-TextCell "D30" "133.93"
Set-TextCell "E30" "  +1.84%  "

# Row 31 - ImmutableX
Set-TextCell "D31" "1.209"
Set-TextCell "E31" "  +2.18%  "

# Row 32 - Stellar
Set-TextCell "D32" "0.1071"
Set-TextCell "E32" "  +0.27%  "

# Row 33 - ARBITRUM
Set-TextCell "D33" "1.651"
Set-TextCell "E33" "  -0.43%  "

# Row 34 - Filecoin
Set-TextCell "D34" "6.362"
Set-TextCell "E34" "  +3.28%  "

# Row 35 - HuobiToken
Set-TextCell "D35" "3.957"
Set-TextCell "E35" "  +1.52%  "

# Row 36 - FraxShare
Set-TextCell "D36" "10.53"
Set-TextCell "E36" "  +4.70%  "

# Row 37 - InternetComputer(DFINITY)
Set-TextCell "D37" "5.923"
Set-TextCell "E37" "  +8.35%  "

# Row 38 - VeChain
Set-TextCell "D38" "0.02624"
Set-TextCell "E38" "  +2.14%  "

# Row 39 - Hedera
Set-TextCell "D39" "0.06825"
Set-TextCell "E39" "  +0.00%  "

# Row 40 - Algorand
Set-TextCell "D40" "0.2328"
Set-TextCell "E40" "  +3.17%  "

# Row 41 - Aptos
Set-TextCell "D41" "12.60"
Set-TextCell "E41" "  -0.49%  "

# Row 42 - TheSandbox
Set-TextCell "D42" "0.6871"
Set-TextCell "E42" "  -0.24%  "

# Row 43 - TrustWalletToken (price unchanged)
Set-TextCell "E43" "  +0.09%  "

# Row 44 - EnergySwap (price unchanged)
Set-TextCell "E44" "  +7.55%  "

# Row 45 - Decentraland
Set-TextCell "D45" "0.6451"
Set-TextCell "E45" "  +1.58%  "

# Row 46 - NEARProtocol
Set-TextCell "D46" "2.315"
Set-TextCell "E46" "  +5.32%  "

# Row 47 - was PancakeSwap, now BabyDogeCoin (ranking swap with row 48)
Set-TextCell "B47" "BabyDogeCoin"
Set-TextCell "C47" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D47" "0.00000000367"
Set-TextCell "E47" "  +23.65%  "

# Row 48 - was BabyDogeCoin, now PancakeSwap (ranking swap with row 47)
Set-TextCell "B48" "PancakeSwap"
Set-TextCell "C48" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D48" "3.689"
Set-TextCell "E48" "  +1.91%  "

# Row 49 - EOS
Set-TextCell "D49" "1.254"
Set-TextCell "E49" "  +0.87%  "

# Row 50 - Aave
Set-TextCell "D50" "83.42"
Set-TextCell "E50" "  +1.92%  "

# Row 51 - WOONetwork
Set-TextCell "D51" "0.3358"
Set-TextCell "E51" "  +12.54%  "
